# issue #5: stock data output to json file
#
# The 股票 (stock) sheet gains a new "property_category" column (value
# "stock" for every row), inserted right after the existing "total"
# column and before the "date" column. The columns that used to sit at
# H/I/J (date / legislator_name / legislator_id) shift one column to
# the right, to I/J/K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # 股票 (stock) sheet

# Insert a new blank column at H, shifting existing H:J -> I:K.
$ws.Columns.Item(8).Insert(-4161)  # -4161 == xlShiftToRight

# Header for the new column.
$ws.Range("H1").Value = "property_category"

# Value for every data row in the new column.
$lastRow = $ws.Cells(($ws.Rows.Count), 1).End(-4162).Row  # -4162 == xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("H" + $r).Value = "stock"
}
